$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 250, shifting existing rows 250-329 down to 251-330.
$ws.Rows(250).Insert()

# Populate the newly inserted row 250 with the new record's data.
$ws.Range("A250").Value = 7
$ws.Range("B250").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C250").Value = "Ñuble"
$ws.Range("D250").Value = 45120
$ws.Range("D250").NumberFormat = $ws.Range("D251").NumberFormat
$ws.Range("E250").Value = 16
$ws.Range("F250").Value = 100112032
$ws.Range("G250").Value = "Zapallo italiano"
$ws.Range("H250").Value = "Sin especificar"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 80
$ws.Range("K250").Value = 14000
$ws.Range("L250").Value = 14000
$ws.Range("M250").Value = 14000
$ws.Range("N250").Value = '$/caja 50 unidades'
$ws.Range("O250").Value = "Región de Arica y Parinacota"
$ws.Range("P250").Value = 280
$ws.Range("Q250").Value = 50
$ws.Range("R250").Value = "Hortaliza"
